{"js": "// Replace the two-digit x two-digit multiplication problems in the\n// worksheet table with a new set of problems, cell by cell (positional\n// match \u2014 two of the old problems repeat the same text with different\n// replacements, so a plain global find/replace would be ambiguous).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Rows with visible math problems: 0, 4, 9, 14, 19 (5 columns each).\nconst replacements = [\n  [\"29\u00d770=\", \"80\u00d798=\", \"25\u00d720=\", \"17\u00d787=\", \"32\u00d788=\"],\n  [\"12\u00d726=\", \"71\u00d784=\", \"44\u00d730=\", \"84\u00d731=\", \"20\u00d727=\"],\n  [\"69\u00d779=\", \"45\u00d796=\", \"18\u00d721=\", \"96\u00d757=\", \"85\u00d737=\"],\n  [\"62\u00d752=\", \"58\u00d772=\", \"98\u00d745=\", \"45\u00d785=\", \"73\u00d746=\"],\n  [\"60\u00d765=\", \"97\u00d725=\", \"84\u00d736=\", \"21\u00d750=\", \"46\u00d772=\"],\n];\nconst rowIndexes = [0, 4, 9, 14, 19];\n\nfor (let r = 0; r < rowIndexes.length; r++) {\n  const rowIdx = rowIndexes[r];\n  for (let c = 0; c < 5; c++) {\n    table.getCell(rowIdx, c).value = replacements[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the two-digit x two-digit multiplication problems in the\n# worksheet table with a new set of problems, cell by cell (positional\n# match -- two of the old problems repeat the same text with different\n# replacements, so a plain global find/replace would be ambiguous).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row 1 (table row 1 / 1-indexed)\n$t.Cell(1, 1).Range.Text = \"29\u00d770=\"\n$t.Cell(1, 2).Range.Text = \"80\u00d798=\"\n$t.Cell(1, 3).Range.Text = \"25\u00d720=\"\n$t.Cell(1, 4).Range.Text = \"17\u00d787=\"\n$t.Cell(1, 5).Range.Text = \"32\u00d788=\"\n\n# Row 5\n$t.Cell(5, 1).Range.Text = \"12\u00d726=\"\n$t.Cell(5, 2).Range.Text = \"71\u00d784=\"\n$t.Cell(5, 3).Range.Text = \"44\u00d730=\"\n$t.Cell(5, 4).Range.Text = \"84\u00d731=\"\n$t.Cell(5, 5).Range.Text = \"20\u00d727=\"\n\n# Row 10\n$t.Cell(10, 1).Range.Text = \"69\u00d779=\"\n$t.Cell(10, 2).Range.Text = \"45\u00d796=\"\n$t.Cell(10, 3).Range.Text = \"18\u00d721=\"\n$t.Cell(10, 4).Range.Text = \"96\u00d757=\"\n$t.Cell(10, 5).Range.Text = \"85\u00d737=\"\n\n# Row 15\n$t.Cell(15, 1).Range.Text = \"62\u00d752=\"\n$t.Cell(15, 2).Range.Text = \"58\u00d772=\"\n$t.Cell(15, 3).Range.Text = \"98\u00d745=\"\n$t.Cell(15, 4).Range.Text = \"45\u00d785=\"\n$t.Cell(15, 5).Range.Text = \"73\u00d746=\"\n\n# Row 20\n$t.Cell(20, 1).Range.Text = \"60\u00d765=\"\n$t.Cell(20, 2).Range.Text = \"97\u00d725=\"\n$t.Cell(20, 3).Range.Text = \"84\u00d736=\"\n$t.Cell(20, 4).Range.Text = \"21\u00d750=\"\n$t.Cell(20, 5).Range.Text = \"46\u00d772=\"\n"}
